# Updated symbol list on Wed Jan  4 09:30:25 UTC 2023 with GitHub Actions
# Applies refreshed price/volume figures (and a symbol-list reshuffle in
# rows 9-14) to the cryptos worksheet, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''255.27'
$ws.Range("E2").Value = '''3.70%'
$ws.Range("D3").Value = '''28.11'
$ws.Range("E3").Value = '''-5.58%'
$ws.Range("D4").Value = '''5.255'
$ws.Range("E4").Value = '''1.97%'
$ws.Range("D5").Value = '''0.05851'
$ws.Range("E5").Value = '''1.49%'
$ws.Range("D6").Value = '''6.704'
$ws.Range("E6").Value = '''0.77%'
$ws.Range("E7").Value = '''2.24%'
$ws.Range("D8").Value = '''1.039'
$ws.Range("E8").Value = '''21.73%'
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = '''0.01062'
$ws.Range("E9").Value = '''1,682.70%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1412'
$ws.Range("E10").Value = '''2.34%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.07112'
$ws.Range("E11").Value = '''0.50%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03187'
$ws.Range("E12").Value = '''-1.84%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09225'
$ws.Range("E13").Value = '''-1.54%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001549'
$ws.Range("E14").Value = '''0.53%'
$ws.Range("D15").Value = '''0.005827'
$ws.Range("E15").Value = '''-2.35%'
$ws.Range("E16").Value = '''-0.65%'
$ws.Range("E17").Value = '''-0.17%'
$ws.Range("D18").Value = '''2.202'
$ws.Range("E18").Value = '''-0.94%'
$ws.Range("D19").Value = '''0.3181'
$ws.Range("E19").Value = '''0.53%'
$ws.Range("D20").Value = '''0.03462'
$ws.Range("E20").Value = '''2.73%'
$ws.Range("E21").Value = '''0.03%'
$ws.Range("D22").Value = '''3.566'
$ws.Range("E22").Value = '''1.93%'
$ws.Range("D23").Value = '''0.04149'
$ws.Range("E23").Value = '''0.64%'
$ws.Range("D24").Value = '''0.1347'
$ws.Range("E24").Value = '''-4.43%'
$ws.Range("D25").Value = '''0.001222'
$ws.Range("E25").Value = '''-0.48%'
$ws.Range("E27").Value = '''0.06%'
$ws.Range("D28").Value = '''0.00008002'
$ws.Range("E28").Value = '''-44.75%'
$ws.Range("D40").Value = '''0.03822'
$ws.Range("E40").Value = '''1.99%'
$ws.Range("D41").Value = '''0.005751'
$ws.Range("E41").Value = '''0.94%'
$ws.Range("E42").Value = '''3.06%'
$ws.Range("D43").Value = '''0.002300'
$ws.Range("E43").Value = '''0.06%'
$ws.Range("D44").Value = '''0.009684'
$ws.Range("E44").Value = '''13.35%'
$ws.Range("D45").Value = '''0.00005231'
$ws.Range("E45").Value = '''-4.89%'
$ws.Range("E46").Value = '''0.06%'
$ws.Range("E47").Value = '''31.07%'
$ws.Range("E48").Value = '''-3.46%'
$ws.Range("E49").Value = '''0.06%'
$ws.Range("E50").Value = '''0.06%'

Write-Output "Applied crypto data refresh"
